# FAST_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer banner
#  - refresh the Weight (D) / Percent Change (E) columns for rows 2-10
#
# The sheet ships protected (sheetProtection, password "D382"), so we have
# to unprotect before writing and restore protection with the same
# password afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect("D382")

# --- Update the "as of" date inside the confidential disclaimer text (A13) ---
$disclaimerCell = $ws.Range("A13")
$disclaimerCell.Value = $disclaimerCell.Text.Replace("2021-06-14", "2021-07-07")
# Re-run the row's autofit so the edit doesn't leave a stray explicit row
# height behind (writing a wrapped/multi-line string otherwise stamps the
# row with a custom height in the saved XML).
$ws.Rows(13).EntireRow.AutoFit()

# --- Refresh Weight (D) / Percent Change (E) for each holding row ---
$weights = @{
    2  = 0.1386701192887695
    3  = 0.1085024521476936
    4  = 0.1124073359632902
    5  = 0.1187233503483783
    6  = 0.1209905445892478
    7  = 0.1425887188333804
    8  = 0.1301565928722872
    9  = 0.127960885956953
    10 = 1
}

foreach ($row in $weights.Keys) {
    $ws.Cells.Item($row, 4).Value = $weights[$row]
    $ws.Cells.Item($row, 5).Value = 0
}

# --- Restore the original sheet protection ---
$ws.Protect("D382")
